$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target OOXML width 15.42578125; the COM layer
# quantizes ColumnWidth to 1/6-character increments, so 14.6 is the input
# that lands closest on the achievable grid, 15.5)
$ws.Columns.Item(1).ColumnWidth = 14.6

# Update values in A1:A33
$values = @(
    0.21342162492258865,
    -0.0059999999655673264,
    -0.0039999999691193722,
    -0.0079999999436939362,
    -0.0029999999680452305,
    -0.0019999999656956646,
    -0.0099999999214390733,
    -0.0099999999202751155,
    -0.0019999999640503141,
    -0.0019999999637718702,
    -0.0029999999583401049,
    -0.0034999999555607175,
    -0.0034999999562268513,
    -0.0079999999318669524,
    0.020158417041494125,
    -0.0019999999654389811,
    -0.0019999999649336075,
    -0.0039999999536757258,
    -0.0039999999760631511,
    -0.029596668024741035,
    -0.0039999999702109434,
    -0.0039999999700262023,
    -0.0049999999617593716,
    -0.01999999987556933,
    -0.019999999873957286,
    -0.00249999996194461,
    -0.0076682339673705968,
    -0.0019999999572073435,
    -0.0069999999252559064,
    -0.059999999630198086,
    -0.0069999999244281241,
    0.037966594268054266,
    -0.018271729722421526
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
